# Add a "PageObject" column (Object Repository) between the existing
# "Description" and "Action_Keyword" columns on the TestSteps sheet,
# and populate it with the new page-object locator names for the
# login-related test steps.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D; this shifts the old column D (Action_Keyword)
# to column E and carries over cell styles/formatting automatically.
$ws.Columns("D:D").Insert()

# Header for the newly inserted column.
$ws.Range("D1").Value = "PageObject"

# Populate the PageObject values for the relevant test steps.
# Rows 2, 3 and 10 (Open Browser / Navigate / Close Browser) intentionally
# have no page object, so they are left blank.
$ws.Range("D4").Value = "login_txtBox_empName"
$ws.Range("D5").Value = "login_txtBox_empNumber"
$ws.Range("D6").Value = "login_txtBox_password"
$ws.Range("D7").Value = "login_txtBox_CnfrmPassword"
$ws.Range("D8").Value = "login_btn_submit"
$ws.Range("D9").Value = "login_div_result"

# Re-size the Description column (now narrower since it only holds the
# short description) and the new PageObject column.
$ws.Columns("C:C").ColumnWidth = 26.75
$ws.Columns("D:D").ColumnWidth = 29.27

# Update the active selection to match the new layout.
$ws.Range("D9").Select()
